$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("J2:J16").Formula = "=G2*0.08"
$ws.Columns("L:L").Delete()
$ws.Columns("B:B").ColumnWidth = 10.666666666666666
$ws.Columns("J:J").ColumnWidth = 10.498697916666666
$ws.Columns("K:K").ColumnWidth = 9.330729166666666
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").HorizontalAlignment = -4108
$ws.Range("L21").Select()
